# [Common] Update email, formula
# - Replace the stale "=A2+B2+AC1" placeholder formula in column K (rows 2-18)
#   with the real profit/loss check: =IF(E2-(E2*0.15+J2)<0, TRUE, FALSE)
# - Correct the eBay cost (J8) picked up for the "Handle Bridge Kit" row
# - Fix a product-title typo in H14 ("Self-Cleaning" -> "Self Cleaning")
# - Move the active selection to D23

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the profit-check formula for every data row (2 through 18).
for ($r = 2; $r -le 18; $r++) {
    $ws.Range("K$r").Formula = "=IF(E2-(E2*0.15+J2)<0, TRUE, FALSE)"
}

# Corrected Amazon price for the "Handle Bridge Kit" listing (row 8).
$ws.Range("J8").Value = 13.88

# Typo fix in the Amazon item title for row 14.
$ws.Range("H14").Value = "PetSafe Simply Clean Self Cleaning Cat Litter Box, Automatic Litter Box, Works with Clumping Cat Litter"

# Restore the sheet's active selection.
$ws.Range("D23").Select()
